$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Check the "Check Box 14" control bound to cell J18 (ctrlProp5), which marks
# checklist item 5 (row 18) as done. Setting the linked cell to TRUE mirrors
# what happens when the checkbox is clicked/checked.
$ws.Range("J18").Value = $true

$excel.CalculateFullRebuild()
$wb.Save()
